# The "sold car" / 5000 / 7/7/2025 entry (row 8) is removed from the
# Incomes table. Deleting the entire row shifts the two rows below it
# (abhayawas@gmail.com and john@example.com) up by one, which is exactly
# what the target sheet looks like (now ending at row 9 instead of row 10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Delete()
